$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-28 changes from serial date 45440 (2024-05-28)
# to serial date 45441 (2024-05-29).
$newDate = (Get-Date -Year 2024 -Month 5 -Day 29 -Hour 0 -Minute 0 -Second 0).Date

for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
